# WIP + specify binding by hand
#
# Formal Ampersand "Braga" workbook: replace the automatically-derived
# identity-composition terms (t9="I[A];r", t10="I[A];s", t11="I[B];t" and
# their downstream Binary Terms / Compositions rows) with a hand-specified
# `bind` relationship: the base relation Terms (t1=r, t2=s, t3=t) now each
# point, via a new "bind" column, straight at typed relation Terms
# (t6=r[A*C], t7=s[A*B], t8=t[B*C]) instead of going through I[x] atoms.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Relations sheet: t1/t2/t3 -> t6/t7/t8 (same relations, renumbered Terms)
# ---------------------------------------------------------------------
$wsRelations = $wb.Worksheets.Item("Relations")
$wsRelations.Range("A3").Value = "t6"
$wsRelations.Range("A4").Value = "t7"
$wsRelations.Range("A5").Value = "t8"

# ---------------------------------------------------------------------
# Terms sheet: add "bind" column; re-purpose t6/t7/t8; drop t9/t10/t11;
# shift the I[A]/I[B]/I[C] identity-term rows up to fill the gap.
# ---------------------------------------------------------------------
$wsTerms = $wb.Worksheets.Item("Terms")

$wsTerms.Range("B9").Value = "r[A*C]"
$wsTerms.Range("B10").Value = "s[A*B]"
$wsTerms.Range("B11").Value = "t[B*C]"

$wsTerms.Range("C1").Value = "bind"
$wsTerms.Range("C2").Value = "Relation"
$wsTerms.Range("C3").Value = "t6"
$wsTerms.Range("C4").Value = "t7"
$wsTerms.Range("C5").Value = "t8"

$wsTerms.Range("A13:B15").ClearContents()
$wsTerms.Range("A17:B19").Cut($wsTerms.Range("A13:B15"))

# ---------------------------------------------------------------------
# Binary Terms sheet: drop the rows built on top of t6..t11
# ---------------------------------------------------------------------
$wsBinaryTerms = $wb.Worksheets.Item("Binary Terms")
$wsBinaryTerms.Range("A9:D15").ClearContents()

# ---------------------------------------------------------------------
# Compositions sheet: drop t6..t11 rows, keep only t5
# ---------------------------------------------------------------------
$wsCompositions = $wb.Worksheets.Item("Compositions")
$wsCompositions.Range("A9:A15").ClearContents()

# ---------------------------------------------------------------------
# Cursor / selection bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------
$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Activate()
$wsConcepts.Range("D18").Select()

$wsBinaryTerms.Activate()
$wsBinaryTerms.Rows("8:15").Select()

$wsCompositions.Activate()
$wsCompositions.Rows("8:15").Select()

$wsTerms.Activate()
$wsTerms.Range("C2").Select()

Write-Output "Braga.xlsx: bind column + manual binding applied"
